$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Lao động"
$ws.Range("B2").Value = "#000080"
$ws.Range("B2").Select() | Out-Null
